$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "timestamp" column (O) for every data row (2-31): the scrape
#    was re-run later the same day, 07:13:58 -> 21:00:37.
# ---------------------------------------------------------------------------
$ws.Range("O2:O31").Value = "2022-08-31 21:00:37"

# ---------------------------------------------------------------------------
# 2) Rows 25/26 and rows 27/28 swap their entire product records (A..N).
#    Some of the swapped-in text looks like a plain number ("3874909",
#    "3.65", "0.12", ...); writing it via .Value would let Excel silently
#    coerce it to a numeric cell, which does not match the source data
#    (plain text). Force those columns to Text before writing, then put the
#    style back to Normal so no stray formatting is left behind.
# ---------------------------------------------------------------------------
$textCols = @("A", "H", "K")
foreach ($col in $textCols) {
    $ws.Range(($col + "25:" + $col + "28")).NumberFormat = "@"
}

# Note: column E ("ratingAmount") is identical before/after the swap in
# every one of these four rows (25<->26 both had 1, 27<->28 both were
# blank), so it is intentionally left untouched to avoid disturbing its
# existing cell representation.
function Set-Row25to28($rowNum, $values) {
    $ws.Range("A" + $rowNum).Value = $values[0]
    $ws.Range("B" + $rowNum).Value = $values[1]
    $ws.Range("C" + $rowNum).Value = $values[2]
    $ws.Range("D" + $rowNum).Value = $values[3]
    $ws.Range("F" + $rowNum).Value = $values[4]
    $ws.Range("G" + $rowNum).Value = $values[5]
    $ws.Range("H" + $rowNum).Value = $values[6]
    $ws.Range("I" + $rowNum).Value = $values[7]
    $ws.Range("J" + $rowNum).Value = $values[8]
    $ws.Range("K" + $rowNum).Value = $values[9]
    $ws.Range("L" + $rowNum).Value = $values[10]
    $ws.Range("M" + $rowNum).Value = $values[11]
    $ws.Range("N" + $rowNum).Value = $values[12]
}

# --- new row 25 (was row 26: "Tempo Bamboo Eco") ---
Set-Row25to28 25 @(
    "6868354",
    "Tempo Bamboo Eco",
    "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco/p/6868354",
    "12ST",
    3,
    "Tempo",
    "3.95",
    "0.33/1ST",
    "Preis pro 1 Stück",
    "0.33",
    "1ST",
    "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']",
    "Tempo Bamboo Eco 3.95 Schweizer Franken"
)

# --- new row 26 (was row 25: "Oecoplan Papiertaschentücher ...") ---
Set-Row25to28 26 @(
    "3874909",
    "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück",
    "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-papiertaschentuecher-special-edition-calendula-30x10-stueck/p/3874909",
    "30ST",
    5,
    "Coop",
    "3.65",
    "0.12/1ST",
    "Preis pro 1 Stück",
    "0.12",
    "1ST",
    "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']",
    "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück 20% Aktion 3.65 Schweizer Franken statt 4.60 Schweizer Franken"
)

# --- new row 27 (was row 28: "Tela Toilettenpapier Futura ...") ---
Set-Row25to28 27 @(
    "6996129",
    "Tela Toilettenpapier Futura 3-lagig 9 Rollen",
    "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tela-toilettenpapier-futura-3-lagig-9-rollen/p/6996129",
    "9Rol",
    0,
    "Tela",
    "8.40",
    "0.93/1Rol",
    "Preis pro 1 Rolle",
    "0.93",
    "1Rol",
    "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']",
    "Tela Toilettenpapier Futura 3-lagig 9 Rollen 8.40 Schweizer Franken"
)

# --- new row 28 (was row 27: "Tela Lux 4-lagig 6 Rollen") ---
Set-Row25to28 28 @(
    "6996029",
    "Tela Lux 4-lagig 6 Rollen",
    "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tela-lux-4-lagig-6-rollen/p/6996029",
    "6Rol",
    0,
    "Tela",
    "5.60",
    "0.93/1Rol",
    "Preis pro 1 Rolle",
    "0.93",
    "1Rol",
    "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']",
    "Tela Lux 4-lagig 6 Rollen 5.60 Schweizer Franken"
)

foreach ($col in $textCols) {
    $ws.Range(($col + "25:" + $col + "28")).Style = "Normal"
}
